$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new log row (row 32) below the last existing row (31).
# Copy the formatting of the previous row so the new row matches the
# existing "data row" style (centered alignment, same font/fill/border).
$ws.Range("A31:H31").Copy()
$ws.Range("A32:H32").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A32").Value = "2025-08-19 09:40:02 UTC"
$ws.Range("B32").Value = "2025-08-19 15:10:02 IST"
$ws.Range("C32").Value = "SKIPPED"
$ws.Range("D32").Value = "No change in PDF. Skipping download & Excel update."
$ws.Range("E32").Value = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$ws.Range("G32").Value = 0
